# Applies the "Added some projects with images and description" commit:
#  - fixes the Cave Cart (D2) and This Website (D3) descriptions
#  - adds 8 new project rows (4-11) with category/name/date/text/image/icon/link
#  - wraps text + taller rows for the two long "text" cells that got a style
#  - leaves the existing G2 hyperlink (unity url) untouched

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (Cave Cart): corrected / expanded description ---
$ws.Range("D2").Value = 'In 2019 I attended a Danish folk high school, where I took the game design course. We were a group of three people on the course. I made several games with my group. Cave Cart was the final project, which took a month to make. My main responsibility on the project was programming using Unity and C#. The game is similar to temple run. The player is driving a mine-cart in a tunnel with three tracks. To avoid obstacles like holes and explosive barrels the player switches tracks. The player gets points according to how long they survive and gets extra points by collecting pickups. The maze is procedurally generated, and the game is endless. The game also contains power-ups, a start screen, setting and credit scene.'
$ws.Range("D2").WrapText = $true
$ws.Rows.Item(2).RowHeight = 16.2

# --- Row 3 (This Website): real description instead of placeholder ---
$ws.Range("D3").Value = 'During my summer vacation after 6.th semester I spent about a month making this personal website to show off my coding projects. I have only ever worked a little with web development, which was during high school. Therefore I had to learn HTML, CSS and JavaScript from scratch, which took some time. I wanted to avoid dublicate code on the website, so I use JavaScript to load the project info from a CSV. So the code in the projects grid and in the popup window is generated by JS. That makes the website easier to maintain and easier to add new projects. The design of the website is my own. I didn''t copy the design or use a temlate. By clicking the GitHub logo above you can see the project in GitHub.'

# --- New project rows 4-11 ---
# Row 4: Dream Game
$ws.Range("A4").Value = 'Dream Game'
$ws.Range("B4").Value = 'unity'
$ws.Range("C4").Value = '01/07/2023 - 01/09/2024'
$ws.Range("D4").Value = 'I went to a danish folk highschool where I took the game design course. One of the first games I made with my group was a 2D platformer called dream game. It’s a simple platformer with a surreal dream-like theme. It was a small project and only took about a week to make. The game can be played online by clicking the Unity icon.'
$ws.Range("E4").Value = 'Background_Image_HOME_edited.jpg'
$ws.Range("F4").Value = 'unity'
$ws.Range("G4").Value = 'https://github.com'

# Row 5: Mind Hive
$ws.Range("A5").Value = 'Mind Hive'
$ws.Range("B5").Value = 'unity'
$ws.Range("C5").Value = '01/07/2023 - 01/09/2025'
$ws.Range("D5").Value = 'This currently unfinished mobile game is a project my friend and I decided to make in our free time. '
$ws.Range("E5").Value = 'Background_Image_HOME_edited.jpg'
$ws.Range("F5").Value = 'unity'
$ws.Range("G5").Value = 'https://github.com'

# Row 6: Flip Maze
$ws.Range("A6").Value = 'Flip Maze'
$ws.Range("B6").Value = 'unity'
$ws.Range("C6").Value = '01/07/2023 - 01/09/2026'
$ws.Range("D6").Value = 'When I was a beginner at Unity I decided I wanted to make a simple puzzle game, to learn the engiene. I made this game called Flip Maze where the goal is to get the red cube through the maze and into the yellow goal. The cube can only be moved by rotating the maze. I have made several features to the game, like a script that can generate the mesh for the maze from a txt file defining the level, and a script that can generate random levels, with a specified difficulty.'
$ws.Range("E6").Value = 'Background_Image_HOME_edited.jpg'
$ws.Range("F6").Value = 'unity'
$ws.Range("G6").Value = 'https://github.com'

# Row 7: Hypercube
$ws.Range("A7").Value = 'Hypercube'
$ws.Range("B7").Value = 'unity'
$ws.Range("C7").Value = '01/07/2023 - 01/09/2027'
$ws.Range("D7").Value = 'On my summer vacation after highschool I had an interest for higher dimensions. I had seen games where the player sees either a 3D slice or a 3D projection of a 4D world. Therefore I made this demo where the user can generate a hypercube in ''n'' dimensions and then rotate and translate it. A hypercube is the higher dimensional generalization of a cube. So a 2D hypercube is a square, a 3D hypercube is a cube and a 4D hypercube is a tesseract. What is shown on the screen is a 3D slice of this n-dimensional hypercube.'
$ws.Range("E7").Value = 'Background_Image_HOME_edited.jpg'
$ws.Range("F7").Value = 'unity'
$ws.Range("G7").Value = 'https://github.com'

# Row 8: Offline RL Using Transformers
$ws.Range("A8").Value = 'Offline RL Using Transformers'
$ws.Range("B8").Value = 'AI / ML'
$ws.Range("C8").Value = '01/07/2023 - 01/09/2028'
$ws.Range("D8").Value = 'During our fourth semester me and my group at DTU made a large machine learning project. All students on fourth semester chose a project from a list. We chose a project about offline reinforcement learning using transformers. My main focus on the project was programming, but also writing the report. The goal of the project was to replicate the results of a research paper, and build on the research paper to obtain new results. Offline RL is the task of training a RL model on existing data without being able to collect new data from the environment. This is generally more difficult than online RL because the model cannot sample new data in places where it lacks knowledge. The researchers proposed that a transformer model could be used for offline RL. This is based on the intuition that RL is about modelling a sequence of actions, and transformers are good at modelling sequences. In the paper the problem was to teach an AI how to walk. That AI could be a simulated spider or humanoid. We used MoJuCo to simulate the environment. In the end the transformer AI was able to walk for about 20 seconds before falling, but didn''t seem to generalize much from the data. That means the model performed well when the environment was similar to the data, but not so well when having to exstrapolate from the data.We exstended on the research by'
$ws.Range("E8").Value = 'Background_Image_HOME_edited.jpg'
$ws.Range("F8").Value = 'github'
$ws.Range("G8").Value = 'https://github.com'

# Row 9: Quadratic Equation Finder
$ws.Range("A9").Value = 'Quadratic Equation Finder'
$ws.Range("B9").Value = 'other programming'
$ws.Range("C9").Value = '01/07/2023 - 01/09/2029'
$ws.Range("D9").Value = 'When I started programming I made small C# console projects in visual studio. One of the first projects I made was a simple algorithm to find a quadratic equation from three points. The algorithm uses an iterative approach to find a, b, and c in the equadratic equation from the three points. The reason it doesn''t just solve three equations with three unknowns is that I hadn''t learned about systems of equations when I made the project. It was still a fun project to learn some basic programming skills.'
$ws.Range("E9").Value = 'Background_Image_HOME_edited.jpg'
$ws.Range("F9").Value = 'github'
$ws.Range("G9").Value = 'https://github.com'

# Row 10: Mavis Robot Project
$ws.Range("A10").Value = 'Mavis Robot Project'
$ws.Range("B10").Value = 'AI / ML'
$ws.Range("C10").Value = '01/07/2023 - 01/09/2030'
$ws.Range("D10").Value = 'During my sixth semester me and my group made a pathfinding project for a robot in the course symbolic artificial intelligence. The goal of the project was to utilize different search algorithms to allow a real robot to navigate a maze by itself, and complete tasks. Some of the search algorithms we learned in the course are DFS, BFS, A*, greedy first search and iterative deepening search algorithms.'
$ws.Range("E10").Value = 'Background_Image_HOME_edited.jpg'
$ws.Range("F10").Value = 'github'
$ws.Range("G10").Value = 'https://github.com'

# Row 11: Database Project
$ws.Range("A11").Value = 'Database Project'
$ws.Range("B11").Value = 'other programming'
$ws.Range("C11").Value = '01/07/2023 - 01/09/2031'
$ws.Range("D11").Value = 'On my sixth semester I made a database in MySQL with my group. It is a formel 1 database, that keeps track of races, drivers, tournaments and more. After this course I had a good understanding of how classic databases works and how to maintain them. I learned about SQL queries, how databases works on a fundamental level and database normalization.'
$ws.Range("E11").Value = 'Background_Image_HOME_edited.jpg'
$ws.Range("F11").Value = 'github'
$ws.Range("G11").Value = 'https://github.com'

# Row 8 (Offline RL) text is long too -> wrap + shorter custom row height
$ws.Range("D8").WrapText = $true
$ws.Rows.Item(8).RowHeight = 13.8

# Column A needs to grow to fit the new, longer project names
$ws.Columns.Item(1).AutoFit()

# Matches the saved selection/active cell in the edited workbook
$ws.Range("D3").Select() | Out-Null
